# Apply the "today" function addition to the adate form workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "survey" sheet: fix a stray typo and add the new rows that
#    wire up the "today" / "displayToday" calculations.
# ---------------------------------------------------------------
$wsSurvey = $wb.Worksheets.Item("survey")

# Fix the pre-existing typo: H4 said "Nammename" but should match G4.
$wsSurvey.Cells.Item(4, 8).Value = "Testdato 2"

# Insert a new row 7: an "assign" field that stores today's date into "test".
$wsSurvey.Rows.Item(7).Insert()
$wsSurvey.Cells.Item(7, 4).Value = "assign"
$wsSurvey.Cells.Item(7, 6).Value = "test"
$wsSurvey.Cells.Item(7, 9).Value = "adate.today()"

# Insert a new row 8: a note displaying the "today" calculation.
$wsSurvey.Rows.Item(8).Insert()
$wsSurvey.Cells.Item(8, 4).Value = "note"
$wsSurvey.Cells.Item(8, 7).Value = "today: {{calculates.today}}"
$wsSurvey.Cells.Item(8, 8).Value = "today: {{calculates.today}}"

# After the existing "display" note (now row 11), insert 4 new rows.
$wsSurvey.Range("A12:A15").EntireRow.Insert()

$wsSurvey.Cells.Item(12, 4).Value = "note"
$wsSurvey.Cells.Item(12, 7).Value = "display today: {{calculates.displayToday}}"
$wsSurvey.Cells.Item(12, 8).Value = "display today: {{calculates.displayToday}}"

$wsSurvey.Cells.Item(13, 4).Value = "note"
$wsSurvey.Cells.Item(13, 7).Value = "assign: {{data.test}}"
$wsSurvey.Cells.Item(13, 8).Value = "assign: {{data.test}}"

$wsSurvey.Cells.Item(14, 4).Value = "note"
$wsSurvey.Cells.Item(14, 7).Value = "ADA: {{data.ADA}}"
$wsSurvey.Cells.Item(14, 8).Value = "ADA: {{data.ADA}}"

$wsSurvey.Cells.Item(15, 4).Value = "adate"
$wsSurvey.Cells.Item(15, 6).Value = "test"
$wsSurvey.Cells.Item(15, 7).Value = "adate"
$wsSurvey.Cells.Item(15, 8).Value = "adate"

# ---------------------------------------------------------------
# 2. "calculates" sheet: add the "today" and "displayToday" calcs.
# ---------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("calculates")
$wsCalc.Cells.Item(9, 1).Value = "today"
$wsCalc.Cells.Item(9, 2).Value = "adate.today()"
$wsCalc.Cells.Item(10, 1).Value = "displayToday"
$wsCalc.Cells.Item(10, 2).Value = "adate.display(adate.today())"

# ---------------------------------------------------------------
# 3. "model" sheet: register the new "ADA3" (adate) and "test"
#    (text) session variables.
# ---------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("model")
$wsModel.Cells.Item(4, 1).Value = "ADA3"
$wsModel.Cells.Item(4, 2).Value = "adate"
$wsModel.Cells.Item(4, 3).Value = $true
$wsModel.Cells.Item(5, 1).Value = "test"
$wsModel.Cells.Item(5, 2).Value = "text"
$wsModel.Cells.Item(5, 3).Value = $true

# ---------------------------------------------------------------
# 4. Restore the cursor / selection positions recorded for each
#    sheet, then leave "survey" as the active tab.
# ---------------------------------------------------------------
$wsCalc.Activate()
$wsCalc.Range("B10").Select() | Out-Null

$wsModel.Activate()
$wsModel.Range("D5").Select() | Out-Null

$wsSurvey.Activate()
$wsSurvey.Range("E16").Select() | Out-Null
